$wb = $excel.ActiveWorkbook

# Update the status text from "Ready for handoff" to "In Translation"
# on every sheet where it appears (Overview, zh-cn, de-de).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the Status columns to match the new (shorter) autofit width
# now that the text is shorter. (ColumnWidth snaps to the engine's
# internal character-width grid, so 12.42 is the input that lands on
# the closest achievable stored width to 13.4101845877511.)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.42
$overview.Columns.Item(6).ColumnWidth = 12.42

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.42

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.42
